# Parameters sheet: move the "foi" (Force of infection) row from the top
# of the parameter list (row 4) down to the bottom (row 7), so the
# remaining rows ("recrate", "infdeath", "susdeath") each shift up by one.
#
# End result (rows 4-7 of the Parameters sheet):
#   4: recrate    (was row 5)
#   5: infdeath    (was row 6)
#   6: susdeath    (was row 7)
#   7: foi         (was row 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# 1. Stash a copy of row 4 ("foi") just below the current data block (row 8
#    is unused), carrying both values and formatting with it.
$ws.Range("A4:H4").Copy($ws.Range("A8:H8"))

# 2. Remove the original row 4; rows 5-8 (including our stashed copy) shift
#    up by one, landing the stashed "foi" data on row 7.
$ws.Rows(4).Delete()

# 3. The source row only had data out to column H with nothing in G, but
#    the rectangular copy above stamped an (empty, unstyled) placeholder
#    into G7 - clear it so the row matches the original "no G cell" shape.
$ws.Range("G7").Clear()

# 4. Leave the whole-row selection on the row that was moved, mirroring
#    the interactive cut/move gesture that produced this edit.
$ws.Rows(4).Select()
